$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly data between row-pair (4,5) and row-pair (6,7):
# columns D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), S (Precio $/Kg)

$cols = @("D","M","N","O","P","S")

foreach ($col in $cols) {
    $v4 = $ws.Range(($col + "4")).Value2
    $v5 = $ws.Range(($col + "5")).Value2
    $v6 = $ws.Range(($col + "6")).Value2
    $v7 = $ws.Range(($col + "7")).Value2

    $ws.Range(($col + "4")).Value2 = $v6
    $ws.Range(($col + "5")).Value2 = $v7
    $ws.Range(($col + "6")).Value2 = $v4
    $ws.Range(($col + "7")).Value2 = $v5
}
